# Natmi following Dr Hou advice
# Rewrite the LR-pair data rows (2-10) on the active sheet, including three
# new rows for the "ECs" sending cluster, matching the regenerated NATMI output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("ECs", "Amh", "Egfr", "ECs", 2, 0.6666666666666666, 0.363728, 1.091184, 0.3031196334509865, 0.3031196334509865, 3, 1, 2.08532, 6.25596, 0.01753772176136817, 0.01753772176136816, 0.7584892729599999, 6.82640345664, 0.005316027791871308, 0.005316027791871307)
  ,@("ECs", "Amh", "Egfr", "FAPs", 2, 0.6666666666666666, 0.363728, 1.091184, 0.3031196334509865, 0.3031196334509865, 3, 1, 101.898173, 305.694519, 0.8569724579756384, 0.8569724579756383, 37.063218668944, 333.568968020496, 0.2597651773391664, 0.2597651773391664)
  ,@("ECs", "Amh", "Egfr", "sCs", 2, 0.6666666666666666, 0.363728, 1.091184, 0.3031196334509865, 0.3031196334509865, 3, 1, 14.921347, 44.76404100000001, 0.1254898202629935, 0.1254898202629935, 5.427311701616001, 48.845805314544, 0.03803842831994877, 0.03803842831994877)
  ,@("FAPs", "Amh", "Egfr", "ECs", 3, 1, 0.520663, 1.561989, 0.4339043947991109, 0.4339043947991109, 3, 1, 2.08532, 6.25596, 0.01753772176136817, 0.01753772176136816, 1.08574896716, 9.77174070444, 0.00760969454702165, 0.007609694547021649)
  ,@("FAPs", "Amh", "Egfr", "FAPs", 3, 1, 0.520663, 1.561989, 0.4339043947991109, 0.4339043947991109, 3, 1, 101.898173, 305.694519, 0.8569724579756384, 0.8569724579756383, 53.054608448699, 477.4914760382911, 0.3718441157374258, 0.3718441157374258)
  ,@("FAPs", "Amh", "Egfr", "sCs", 3, 1, 0.520663, 1.561989, 0.4339043947991109, 0.4339043947991109, 3, 1, 14.921347, 44.76404100000001, 0.1254898202629935, 0.1254898202629935, 7.768993293061001, 69.92093963754901, 0.0544505845146634, 0.05445058451466339)
  ,@("sCs", "Amh", "Egfr", "ECs", 3, 1, 0.3155576666666667, 0.9466730000000001, 0.2629759717499027, 0.2629759717499027, 3, 1, 2.08532, 6.25596, 0.01753772176136817, 0.01753772176136816, 0.6580387134533333, 5.922348421080001, 0.004611999422475208, 0.004611999422475207)
  ,@("sCs", "Amh", "Egfr", "FAPs", 3, 1, 0.3155576666666667, 0.9466730000000001, 0.2629759717499027, 0.2629759717499027, 3, 1, 101.898173, 305.694519, 0.8569724579756384, 0.8569724579756383, 32.15474970947633, 289.3927473852871, 0.2253631648990461, 0.2253631648990461)
  ,@("sCs", "Amh", "Egfr", "sCs", 3, 1, 0.3155576666666667, 0.9466730000000001, 0.2629759717499027, 0.2629759717499027, 3, 1, 14.921347, 44.76404100000001, 0.1254898202629935, 0.1254898202629935, 4.708545442843667, 42.37690898559301, 0.03300080742838134, 0.03300080742838134)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $i + 2
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}
